# Expansão das análises automáticas:
# adiciona as colunas L, M e N (apoio_medio, contribuicoes, media_contribuicoes)
# ao resumo por mencoes_religiosidade.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cabeçalhos da nova área (linha 1), usando o mesmo estilo das demais células de cabeçalho
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# Copia a formatação (fonte em negrito, borda, alinhamento) da última célula
# de cabeçalho existente (K1) para as novas células de cabeçalho.
$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Valores para cada linha de dados (linhas 2 a 7)
$ws.Range("L2").Value = 92.8307967260526
$ws.Range("M2").Value = 187904
$ws.Range("N2").Value = 303.5605815831987

$ws.Range("L3").Value = 87.50944228358395
$ws.Range("M3").Value = 75649
$ws.Range("N3").Value = 358.5260663507109

$ws.Range("L4").Value = 88.14455763556944
$ws.Range("M4").Value = 131371
$ws.Range("N4").Value = 136.8447916666667

$ws.Range("L5").Value = 93.84286760867113
$ws.Range("M5").Value = 72275
$ws.Range("N5").Value = 170.8628841607565

$ws.Range("L6").Value = 17.93451009587296
$ws.Range("M6").Value = 1926
$ws.Range("N6").Value = 14.37313432835821

$ws.Range("L7").Value = 30.65635216359388
$ws.Range("M7").Value = 282
$ws.Range("N7").Value = 15.66666666666667
